$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates sourced from the crypto-price refresh diff: each entry carries
# a leading apostrophe so Excel stores the digit-grouped/percent strings as
# literal text (matching the original inlineStr cells) instead of coercing
# look-alike values (e.g. "1.00") into numbers.
$edits = @(
    @{ Ref = "D2"; Value = "'45.308.34" },
    @{ Ref = "E2"; Value = "'  +2.57%  " },
    @{ Ref = "D3"; Value = "'2.425.76" },
    @{ Ref = "E3"; Value = "'  -0.02%  " },
    @{ Ref = "D5"; Value = "'319.48" },
    @{ Ref = "E5"; Value = "'  +3.75%  " },
    @{ Ref = "D6"; Value = "'102.91" },
    @{ Ref = "E6"; Value = "'  +1.83%  " },
    @{ Ref = "D7"; Value = "'0.514" },
    @{ Ref = "E7"; Value = "'  +0.37%  " },
    @{ Ref = "D8"; Value = "'1.00" },
    @{ Ref = "E9"; Value = "'  +5.11%  " },
    @{ Ref = "D10"; Value = "'35.50" },
    @{ Ref = "E10"; Value = "'  +0.54%  " },
    @{ Ref = "E11"; Value = "'  -0.10%  " },
    @{ Ref = "E12"; Value = "'  -1.46%  " },
    @{ Ref = "D13"; Value = "'18.29" },
    @{ Ref = "E13"; Value = "'  -3.56%  " },
    @{ Ref = "E14"; Value = "'  +1.46%  " },
    @{ Ref = "D15"; Value = "'2.805.07" },
    @{ Ref = "E15"; Value = "'  +0.03%  " },
    @{ Ref = "D16"; Value = "'2.413.59" },
    @{ Ref = "E16"; Value = "'  -2.10%  " },
    @{ Ref = "D17"; Value = "'0.846" },
    @{ Ref = "E17"; Value = "'  +1.38%  " },
    @{ Ref = "D18"; Value = "'45.236.87" },
    @{ Ref = "E18"; Value = "'  +2.50%  " },
    @{ Ref = "D19"; Value = "'12.22" },
    @{ Ref = "E19"; Value = "'  -0.53%  " },
    @{ Ref = "E20"; Value = "'  -1.29%  " },
    @{ Ref = "D21"; Value = "'0.0₃0921" },
    @{ Ref = "E21"; Value = "'  +1.75%  " },
    @{ Ref = "D22"; Value = "'69.09" },
    @{ Ref = "E22"; Value = "'  +0.87%  " },
    @{ Ref = "D23"; Value = "'244.58" },
    @{ Ref = "E23"; Value = "'  +1.73%  " },
    @{ Ref = "E24"; Value = "'  -0.67%  " },
    @{ Ref = "E25"; Value = "'  +0.33%  " },
    @{ Ref = "E26"; Value = "'  -0.04%  " },
    @{ Ref = "D27"; Value = "'25.75" },
    @{ Ref = "E27"; Value = "'  +2.31%  " },
    @{ Ref = "D28"; Value = "'2.18" },
    @{ Ref = "E28"; Value = "'  -6.27%  " },
    @{ Ref = "D29"; Value = "'9.63" },
    @{ Ref = "E29"; Value = "'  +0.21%  " },
    @{ Ref = "D30"; Value = "'49.60" },
    @{ Ref = "E30"; Value = "'  +2.78%  " },
    @{ Ref = "D31"; Value = "'32.91" },
    @{ Ref = "E31"; Value = "'  +0.21%  " },
    @{ Ref = "D32"; Value = "'20.19" },
    @{ Ref = "E32"; Value = "'  +8.05%  " },
    @{ Ref = "E33"; Value = "'  +7.50%  " },
    @{ Ref = "D34"; Value = "'5.22" },
    @{ Ref = "E34"; Value = "'  +0.67%  " },
    @{ Ref = "E35"; Value = "'  +0.13%  " },
    @{ Ref = "D36"; Value = "'0.0761" },
    @{ Ref = "E36"; Value = "'  +0.26%  " },
    @{ Ref = "E37"; Value = "'  -1.30%  " },
    @{ Ref = "D38"; Value = "'4.43" },
    @{ Ref = "E38"; Value = "'  -1.00%  " },
    @{ Ref = "D39"; Value = "'127.33" },
    @{ Ref = "E39"; Value = "'  -1.19%  " },
    @{ Ref = "E40"; Value = "'  -0.40%  " },
    @{ Ref = "E41"; Value = "'  +0.73%  " },
    @{ Ref = "D42"; Value = "'2.19" },
    @{ Ref = "E42"; Value = "'  -4.39%  " },
    @{ Ref = "D43"; Value = "'20.66" },
    @{ Ref = "E43"; Value = "'  -2.84%  " },
    @{ Ref = "E44"; Value = "'  +0.80%  " },
    @{ Ref = "D45"; Value = "'1.938.25" },
    @{ Ref = "E45"; Value = "'  -0.74%  " },
    @{ Ref = "E47"; Value = "'  +2.67%  " },
    @{ Ref = "B48"; Value = "'Stacks" },
    @{ Ref = "C48"; Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx" },
    @{ Ref = "D48"; Value = "'1.79" },
    @{ Ref = "E48"; Value = "'  +9.22%  " },
    @{ Ref = "B49"; Value = "'FraxShare" },
    @{ Ref = "C49"; Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" },
    @{ Ref = "D49"; Value = "'9.11" },
    @{ Ref = "E49"; Value = "'  -3.13%  " },
    @{ Ref = "D50"; Value = "'76.58" },
    @{ Ref = "E50"; Value = "'  +3.96%  " },
    @{ Ref = "D51"; Value = "'4.85" },
    @{ Ref = "E51"; Value = "'  +6.37%  " }
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Ref)
    $cell.Value = $edit.Value
    $cell.Style = 'Normal'
}
